$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.242.56'
$ws.Range("E2").Value = '  -2.58%  '
$ws.Range("D3").Value = '1.556.27'
$ws.Range("E3").Value = '  -4.07%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").Value = '''206.11'
$ws.Range("E6").Value = '  -0.21%  '
$ws.Range("E7").Value = '  -5.34%  '
$ws.Range("D8").Value = '''0.0603'
$ws.Range("E8").Value = '  -1.78%  '
$ws.Range("E9").Value = '  -3.41%  '
$ws.Range("D10").Value = '''17.68'
$ws.Range("E10").Value = '  -2.95%  '
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").Value = '1.773.18'
$ws.Range("E12").Value = '  -4.04%  '
$ws.Range("D13").Value = '1.548.94'
$ws.Range("E13").Value = '  -4.56%  '
$ws.Range("E14").Value = '  -4.71%  '
$ws.Range("D15").Value = '''0.499'
$ws.Range("E15").Value = '  -4.48%  '
$ws.Range("D16").Value = '25.195.12'
$ws.Range("E16").Value = '  -2.72%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0₃0706'
$ws.Range("E17").Value = '  -3.41%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '''58.74'
$ws.Range("E18").Value = '  -3.79%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '''184.04'
$ws.Range("E20").Value = '  -4.08%  '
$ws.Range("E21").Value = '  -3.26%  '
$ws.Range("E22").Value = '  -3.46%  '
$ws.Range("D23").Value = '''5.82'
$ws.Range("E23").Value = '  -3.85%  '
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  -4.11%  '
$ws.Range("D26").Value = '''139.12'
$ws.Range("E26").Value = '  -3.16%  '
$ws.Range("E27").Value = '  -5.20%  '
$ws.Range("D28").Value = '''14.73'
$ws.Range("E28").Value = '  -2.56%  '
$ws.Range("D29").Value = '''6.38'
$ws.Range("E29").Value = '  -4.82%  '
$ws.Range("E30").Value = '  -6.45%  '
$ws.Range("E31").Value = '  -4.23%  '
$ws.Range("E32").Value = '  -3.48%  '
$ws.Range("E33").Value = '  -4.34%  '
$ws.Range("E34").Value = '  -3.23%  '
$ws.Range("E35").Value = '  -4.06%  '
$ws.Range("D36").Value = '1.083.20'
$ws.Range("E36").Value = '  -2.88%  '
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.0149'
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '''0.813'
$ws.Range("E39").Value = '  +6.45%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.27'
$ws.Range("E40").Value = '  -6.69%  '
$ws.Range("E41").Value = '  -4.80%  '
$ws.Range("D42").Value = '''0.757'
$ws.Range("E42").Value = '  -10.30%  '
$ws.Range("D43").Value = '''92.53'
$ws.Range("E43").Value = '  -5.47%  '
$ws.Range("D44").Value = '''5.02'
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("D45").Value = '1.687.88'
$ws.Range("E45").Value = '  -4.02%  '
$ws.Range("E46").Value = '  -2.99%  '
$ws.Range("D47").Value = '''52.17'
$ws.Range("E47").Value = '  -3.75%  '
$ws.Range("E48").Value = '  -4.86%  '
$ws.Range("E49").Value = '  -1.98%  '
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("E51").Value = '  -0.38%  '
